# ---------------------------------------------------------------------------
# Rename Sheet1 -> Regression, add 5 new (mostly empty) trailing sheets,
# extend the Regression sheet's login-style table with new columns/rows,
# update a couple of existing values, and wire up the extra hyperlinks.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- sheets -----------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Regression"

$prev = $ws
$newNames = @("Smoke", "Sanity", "Sheet4", "Sheet5", "Sheet6")
foreach ($n in $newNames) {
    $newSheet = $wb.Worksheets.Add($null, $prev)
    $newSheet.Name = $n
    $prev = $newSheet
}

# --- Regression sheet data --------------------------------------------
# NB: the exact order below reproduces the author's original sharedStrings
# insertion order (new strings appear in the table in first-use order).
$ws.Range("C2").Value = "url"
$ws.Range("C3").Value = "https://ui.cogmento.com/"
$ws.Range("A3").Value = "rahulscreencast9892@gmail.com"
$ws.Range("B3").Value = "Ra987456321@"
$ws.Range("D2").Value = "browser"
$ws.Range("D3").Value = "chrome"
$ws.Range("E2").Value = "firstName"
$ws.Range("B2").Value = "password"
$ws.Range("G2").Value = "email"
$ws.Range("H2").Value = "description"
$ws.Range("E3").Value = "Sumeet"
$ws.Range("F3").Value = "Desai"
$ws.Range("G3").Value = "sumeet.desai@gmail.com"
$ws.Range("H3").Value = "Create a follow up activity"
$ws.Range("F2").Value = "lastName"

# carry the existing row-2 (plain bordered) formatting across the new columns
$ws.Range("A2").Copy()
$ws.Range("C2:H2").PasteSpecial(-4122) | Out-Null

# row 3: D3/E3/F3/H3 are plain-bordered like row 2 (note: no comma-joined
# multi-area ranges here -- only the first area of those is honoured)
$ws.Range("A2").Copy()
foreach ($addr in @("D3", "E3", "F3", "H3")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# row 3: B3/C3/G3 become hyperlink-style + bordered, matching A3/D7
foreach ($addr in @("B3", "C3", "G3")) {
    $ws.Range($addr).Style = "Hyperlink"
}
$ws.Range("D7").Copy()
foreach ($addr in @("B3", "C3", "G3")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# --- hyperlinks ----------------------------------------------------------
# Hyperlinks.Add() always re-stamps the target cell with the built-in
# Hyperlink style (losing the thin border A3/D7 already had), so snapshot a
# clean bordered-hyperlink cell on a scratch range first, run every Add(),
# then restore the clean formatting from the snapshot afterwards.
$ws.Range("Z1").Value = "x"
$ws.Range("D7").Copy()
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:rahulscreencast9892@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:test@rahul.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://ui.cogmento.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Ra987456321@") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:sumeet.desai@gmail.com") | Out-Null

$ws.Range("Z1").Copy()
foreach ($addr in @("A3", "D7", "B3", "C3", "G3")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false
$ws.Range("Z1").Clear() | Out-Null

# --- column widths ---------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 24.5703125
$ws.Range("G1:H1").EntireColumn.ColumnWidth = 24.42578125

# --- selections / active sheet ---------------------------------------
$wsSanity = $wb.Worksheets.Item("Sanity")
$wsSanity.Activate()
$wsSanity.Range("E18").Select()

$ws.Activate()
$ws.Range("H2").Select()
